$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header label from "Lambda" to "Slip Ratio"
$ws.Range("C1").Value = "Slip Ratio"

# Update the Lambda multiplier constant used in column G (was 30, now 2)
$ws.Range("G2").Value = 2

# Update the selected range/active cell to match the saved view state
$ws.Range("A3:B5").Select()
